$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1102.2833
$ws.Range("I137").Value = 867.58
$ws.Range("J137").Value = 2275.8
$ws.Range("K137").Value = 2602.74
$ws.Range("L137").Value = 6827.400000000001
$ws.Range("M137").Value = -52.74000000000024
$ws.Range("N137").Value = -11927.4

$ws.Range("H138").Value = 3404.3572
$ws.Range("I138").Value = 1075.3684
$ws.Range("J138").Value = 6170.0312
$ws.Range("K138").Value = 3226.1052
$ws.Range("L138").Value = 18510.0936
$ws.Range("M138").Value = 1913.8948
$ws.Range("N138").Value = -28790.0936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16884.732
$ws.Range("I32").Value = 19605.549
$ws.Range("J32").Value = 8450.200000000001
$ws.Range("K32").Value = 19605.549
$ws.Range("L32").Value = 8450.200000000001
$ws.Range("M32").Value = -19318.549
$ws.Range("N32").Value = -9024.200000000001

$ws.Range("H74").Value = 791.9836
$ws.Range("I74").Value = 665.0784
$ws.Range("K74").Value = 665.0784
$ws.Range("M74").Value = 208.9216

$ws.Range("H77").Value = 791.9836
$ws.Range("I77").Value = 665.0784
$ws.Range("K77").Value = 3325.392
$ws.Range("M77").Value = 1042.608

$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21802

$ws.Range("H132").Value = 2351.2058
$ws.Range("I132").Value = 1434.2727
$ws.Range("J132").Value = 4032.25
$ws.Range("K132").Value = 4302.8181
$ws.Range("L132").Value = 12096.75
$ws.Range("M132").Value = -1772.8181
$ws.Range("N132").Value = -17156.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 178.6875
$ws.Range("I64").Value = 201.5
$ws.Range("J64").Value = 171.08333
$ws.Range("K64").Value = 201.5
$ws.Range("L64").Value = 171.08333
$ws.Range("M64").Value = 23.5
$ws.Range("N64").Value = -621.0833299999999

$ws.Range("H67").Value = 178.6875
$ws.Range("I67").Value = 201.5
$ws.Range("J67").Value = 171.08333
$ws.Range("K67").Value = 201.5
$ws.Range("L67").Value = 171.08333
$ws.Range("M67").Value = 578.5
$ws.Range("N67").Value = -1731.08333

$ws.Range("H134").Value = 1362.8438
$ws.Range("I134").Value = 1290.3273
$ws.Range("J134").Value = 1806
$ws.Range("K134").Value = 3870.9819
$ws.Range("L134").Value = 5418
$ws.Range("M134").Value = -1335.9819
$ws.Range("N134").Value = -10488

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1526.1702
$ws.Range("I31").Value = 1301.1904
$ws.Range("J31").Value = 3416
$ws.Range("K31").Value = 1301.1904
$ws.Range("L31").Value = 3416
$ws.Range("M31").Value = -1006.1904
$ws.Range("N31").Value = -4006

$ws.Range("H34").Value = 1526.1702
$ws.Range("I34").Value = 1301.1904
$ws.Range("J34").Value = 3416
$ws.Range("K34").Value = 1301.1904
$ws.Range("L34").Value = 3416
$ws.Range("M34").Value = -1099.1904
$ws.Range("N34").Value = -3820

$ws.Range("H58").Value = 1602.375
$ws.Range("I58").Value = 1363.9231
$ws.Range("J58").Value = 2635.6667
$ws.Range("K58").Value = 1363.9231
$ws.Range("L58").Value = 2635.6667
$ws.Range("M58").Value = -1160.9231
$ws.Range("N58").Value = -3041.6667

$ws.Range("H132").Value = 1704.8948
$ws.Range("I132").Value = 1514.95
$ws.Range("J132").Value = 1915.9445
$ws.Range("K132").Value = 4544.85
$ws.Range("L132").Value = 5747.833500000001
$ws.Range("M132").Value = -2014.85
$ws.Range("N132").Value = -10807.8335

$ws.Range("H136").Value = 1602.375
$ws.Range("I136").Value = 1363.9231
$ws.Range("J136").Value = 2635.6667
$ws.Range("K136").Value = 4091.7693
$ws.Range("L136").Value = 7907.000100000001
$ws.Range("M136").Value = -1541.7693
$ws.Range("N136").Value = -13007.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 57.846153
$ws.Range("I14").Value = 57.846153
$ws.Range("K14").Value = 173.538459
$ws.Range("M14").Value = -0.5384589999999889

$ws.Range("H76").Value = 4012
$ws.Range("I76").Value = 964.3333
$ws.Range("J76").Value = 5535.8335
$ws.Range("K76").Value = 2892.9999
$ws.Range("L76").Value = 16607.5005
$ws.Range("M76").Value = -2509.9999
$ws.Range("N76").Value = -17373.5005

$ws.Range("H79").Value = 4012
$ws.Range("I79").Value = 964.3333
$ws.Range("J79").Value = 5535.8335
$ws.Range("K79").Value = 2892.9999
$ws.Range("L79").Value = 16607.5005
$ws.Range("M79").Value = -1566.9999
$ws.Range("N79").Value = -19259.5005

$ws.Range("H113").Value = 704
$ws.Range("I113").Value = 676.4286
$ws.Range("J113").Value = 731.5714
$ws.Range("K113").Value = 2029.2858
$ws.Range("L113").Value = 2194.7142
$ws.Range("M113").Value = 140.7142000000001
$ws.Range("N113").Value = -6534.7142

$ws.Range("H131").Value = 52639676
$ws.Range("I131").Value = 15908.286
$ws.Range("J131").Value = 83336870
$ws.Range("K131").Value = 47724.858
$ws.Range("L131").Value = 250010610
$ws.Range("M131").Value = -42684.858
$ws.Range("N131").Value = -250020690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 28000
$ws.Range("J18").Value = 20666.666
$ws.Range("L18").Value = 20666.666
$ws.Range("N18").Value = -21252.666

$ws.Range("H21").Value = 8025.75
$ws.Range("J21").Value = 9234.333000000001
$ws.Range("L21").Value = 9234.333000000001
$ws.Range("N21").Value = -9580.333000000001

$ws.Range("H29").Value = 31250
$ws.Range("J29").Value = 12500
$ws.Range("L29").Value = 12500
$ws.Range("N29").Value = -13080

$ws.Range("H30").Value = 8025.75
$ws.Range("J30").Value = 9234.333000000001
$ws.Range("L30").Value = 9234.333000000001
$ws.Range("M30").Value = -4295
$ws.Range("N30").Value = -9444.333000000001

$ws.Range("H92").Value = 23245.5
$ws.Range("J92").Value = 23245.5
$ws.Range("L92").Value = 23245.5
$ws.Range("N92").Value = -26989.5

$ws.Range("H132").Value = 1706
$ws.Range("I132").Value = 1126.96
$ws.Range("K132").Value = 3380.88
$ws.Range("M132").Value = -850.8800000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 15163.4
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 18704.25
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 18704.25
$ws.Range("M17").Value = -830
$ws.Range("N17").Value = -19044.25

$ws.Range("H23").Value = 5251.75

$ws.Range("H132").Value = 3723.4
$ws.Range("I132").Value = 3260.0605
$ws.Range("J132").Value = 4997.5835
$ws.Range("K132").Value = 9780.181500000001
$ws.Range("L132").Value = 14992.7505
$ws.Range("M132").Value = -7250.181500000001
$ws.Range("N132").Value = -20052.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1916.8334
$ws.Range("I96").Value = 2140.2
$ws.Range("J96").Value = 800
$ws.Range("K96").Value = 2140.2
$ws.Range("L96").Value = 800
$ws.Range("M96").Value = -767.1999999999998
$ws.Range("N96").Value = -3546

$ws.Range("H126").Value = 5442.8335
$ws.Range("I126").Value = 6612.357
$ws.Range("J126").Value = 1349.5
$ws.Range("K126").Value = 19837.071
$ws.Range("L126").Value = 4048.5
$ws.Range("M126").Value = -17367.071
$ws.Range("N126").Value = -8988.5

$ws.Range("H132").Value = 1633.2122
$ws.Range("I132").Value = 1241.2727
$ws.Range("J132").Value = 2417.0908
$ws.Range("K132").Value = 3723.8181
$ws.Range("L132").Value = 7251.2724
$ws.Range("M132").Value = -1193.8181
$ws.Range("N132").Value = -12311.2724
